$d = $word.ActiveDocument

# Mapping of old equation text -> new equation text, applied via Find & Replace
# Using unique exact-match search (no wildcards) across the whole document story.
$pairs = @(
    ,@('79×27=2133', '82×75=6150')
    ,@('90×48=4320', '74×60=4440')
    ,@('66×17=1122', '29×91=2639')
    ,@('96×37=3552', '22×37=814')
    ,@('82×48=3936', '33×97=3201')
    ,@('83×13=1079', '60×51=3060')
    ,@('12×81=972', '28×12=336')
    ,@('91×91=8281', '23×24=552')
    ,@('34×52=1768', '100×25=2500')
    ,@('63×85=5355', '81×80=6480')
    ,@('77×68=5236', '93×47=4371')
    ,@('94×75=7050', '20×43=860')
    ,@('18×29=522', '72×93=6696')
    ,@('79×85=6715', '27×80=2160')
    ,@('86×37=3182', '69×27=1863')
    ,@('40×42=1680', '57×80=4560')
    ,@('26×29=754', '88×60=5280')
    ,@('17×80=1360', '23×91=2093')
    ,@('77×37=2849', '30×55=1650')
    ,@('92×61=5612', '85×90=7650')
    ,@('86×63=5418', '83×23=1909')
    ,@('16×96=1536', '63×67=4221')
    ,@('70×19=1330', '80×64=5120')
    ,@('19×53=1007', '60×98=5880')
    ,@('17×46=782', '53×37=1961')
    ,@('100×65=6500', '98×45=4410')
    ,@('78×95=7410', '57×64=3648')
    ,@('30×81=2430', '50×11=550')
    ,@('43×77=3311', '93×95=8835')
    ,@('67×38=2546', '67×58=3886')
    ,@('33×23=759', '65×97=6305')
    ,@('95×23=2185', '87×77=6699')
    ,@('39×55=2145', '67×63=4221')
    ,@('44×24=1056', '65×38=2470')
    ,@('54×91=4914', '35×96=3360')
    ,@('23×12=276', '65×92=5980')
    ,@('60×20=1200', '63×17=1071')
    ,@('25×91=2275', '22×67=1474')
    ,@('23×55=1265', '11×47=517')
    ,@('32×68=2176', '63×100=6300')
    ,@('36×10=360', '49×20=980')
    ,@('43×98=4214', '21×57=1197')
    ,@('96×17=1632', '83×44=3652')
    ,@('75×43=3225', '85×69=5865')
    ,@('71×49=3479', '70×38=2660')
    ,@('30×91=2730', '37×67=2479')
    ,@('67×19=1273', '18×33=594')
    ,@('19×93=1767', '100×84=8400')
    ,@('94×22=2068', '89×39=3471')
    ,@('32×38=1216', '57×59=3363')
    ,@('96×15=1440', '64×50=3200')
    ,@('37×40=1480', '30×45=1350')
    ,@('14×39=546', '60×27=1620')
    ,@('100×41=4100', '63×75=4725')
    ,@('76×43=3268', '16×21=336')
    ,@('58×98=5684', '38×73=2774')
    ,@('27×87=2349', '18×28=504')
    ,@('23×11=253', '100×97=9700')
    ,@('86×34=2924', '95×35=3325')
    ,@('91×86=7826', '19×90=1710')
    ,@('32×80=2560', '50×70=3500')
    ,@('96×14=1344', '88×51=4488')
    ,@('66×90=5940', '59×21=1239')
    ,@('72×35=2520', '49×46=2254')
    ,@('89×58=5162', '19×41=779')
    ,@('100×40=4000', '37×59=2183')
    ,@('35×85=2975', '22×72=1584')
    ,@('19×31=589', '28×36=1008')
    ,@('36×22=792', '35×20=700')
    ,@('63×54=3402', '48×91=4368')
    ,@('25×62=1550', '16×31=496')
    ,@('81×64=5184', '30×44=1320')
    ,@('58×38=2204', '97×90=8730')
    ,@('48×58=2784', '77×30=2310')
    ,@('65×72=4680', '60×99=5940')
    ,@('58×85=4930', '89×17=1513')
    ,@('62×57=3534', '94×35=3290')
    ,@('26×98=2548', '57×96=5472')
    ,@('93×42=3906', '43×66=2838')
    ,@('56×50=2800', '85×78=6630')
    ,@('83×65=5395', '96×67=6432')
    ,@('12×47=564', '57×39=2223')
    ,@('29×81=2349', '97×20=1940')
    ,@('22×50=1100', '50×19=950')
    ,@('14×66=924', '79×100=7900')
    ,@('53×24=1272', '91×98=8918')
    ,@('90×40=3600', '89×16=1424')
    ,@('73×65=4745', '20×92=1840')
    ,@('21×74=1554', '35×16=560')
    ,@('58×79=4582', '94×35=3290')
    ,@('91×73=6643', '67×80=5360')
    ,@('52×71=3692', '59×64=3776')
    ,@('17×74=1258', '32×26=832')
    ,@('51×22=1122', '70×66=4620')
    ,@('54×53=2862', '57×72=4104')
    ,@('30×60=1800', '94×23=2162')
    ,@('44×34=1496', '71×23=1633')
    ,@('24×32=768', '35×28=980')
    ,@('63×60=3780', '98×78=7644')
    ,@('12×83=996', '100×43=4300')
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
